$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on Price/Volume columns so values like "1.005" or "20.747.75" stay as literal text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '20.671.84'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').Value = '1.507.57'
$ws.Range('E3').Value = '  +3.93%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '0.9623'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').Value = '278.84'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('D7').Value = '0.3566'
$ws.Range('E7').Value = '  -1.99%  '
$ws.Range('D8').Value = '0.3110'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').Value = '1.099'
$ws.Range('E9').Value = '  +6.77%  '
$ws.Range('D10').Value = '39.57'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = '0.06713'
$ws.Range('E11').Value = '  +2.97%  '
$ws.Range('D12').Value = '0.9993'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '18.46'
$ws.Range('E13').Value = '  +4.65%  '
$ws.Range('D14').Value = '5.559'
$ws.Range('E14').Value = '  +3.78%  '
$ws.Range('D15').Value = '6.232'
$ws.Range('E15').Value = '  +2.59%  '
$ws.Range('D16').Value = '0.9628'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').Value = '0.00001023'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '1.501.60'
$ws.Range('E18').Value = '  +3.84%  '
$ws.Range('D19').Value = '0.06033'
$ws.Range('E19').Value = '  +5.80%  '
$ws.Range('D20').Value = '69.90'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '5.569'
$ws.Range('E21').Value = '  +3.27%  '
$ws.Range('D22').Value = '14.72'
$ws.Range('E22').Value = '  +2.80%  '
$ws.Range('D23').Value = '11.26'
$ws.Range('E23').Value = '  +4.25%  '
$ws.Range('E24').Value = '  +3.08%  '
$ws.Range('D25').Value = '20.732.28'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Value = '146.70'
$ws.Range('E26').Value = '  +4.23%  '
$ws.Range('D27').Value = '2.138'
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('D28').Value = '17.46'
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('D29').Value = '1.664.10'
$ws.Range('E29').Value = '  +4.04%  '
$ws.Range('D30').Value = '115.93'
$ws.Range('E30').Value = '  +4.27%  '
$ws.Range('D31').Value = '3.979'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').Value = '5.039'
$ws.Range('E32').Value = '  +4.68%  '
$ws.Range('D33').Value = '0.8236'
$ws.Range('E33').Value = '  +5.06%  '
$ws.Range('D34').Value = '0.07961'
$ws.Range('E34').Value = '  +3.00%  '
$ws.Range('D35').Value = '1.197'
$ws.Range('E35').Value = '  +7.19%  '
$ws.Range('D36').Value = '1.446'
$ws.Range('E36').Value = '  -2.59%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '4.813'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05754'
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('D39').Value = '0.02048'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('D40').Value = '10.49'
$ws.Range('E40').Value = '  +2.93%  '
$ws.Range('D41').Value = '0.9625'
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('D42').Value = '7.532'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').Value = '0.1877'
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('D44').Value = '0.5292'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '12.31'
$ws.Range('E45').Value = '  +3.64%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.537'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('D47').Value = '119.92'
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('D48').Value = '0.5257'
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('D49').Value = '1.849'
$ws.Range('E49').Value = '  +6.57%  '
$ws.Range('D50').Value = '0.06480'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').Value = '0.9831'
$ws.Range('E51').Value = '  -0.14%  '
